# Corrects field-order / value drift produced by re-running the reflection-
# based "Feign Relations Finder" analyzer against the mall-common module
# (non-deterministic member ordering) and records the associated constant
# swaps, matching commit "Fixed Feign Relations Finder and added 7 projects."

$wb = $excel.ActiveWorkbook

$wsClassFields = $wb.Worksheets.Item("classFields")
$wsClassFields.Range("D4").Value = 'java.util.List'
$wsClassFields.Range("D5").Value = 'java.lang.String'
$wsClassFields.Range("D6").Value = 'java.lang.Long'
$wsClassFields.Range("D7").Value = 'java.lang.String'
$wsClassFields.Range("D9").Value = 'java.lang.Object'
$wsClassFields.Range("D10").Value = 'java.lang.String'
$wsClassFields.Range("B12").Value = 'roles'
$wsClassFields.Range("D13").Value = 'java.lang.String'
$wsClassFields.Range("D14").Value = 'java.lang.Integer'
$wsClassFields.Range("D16").Value = 'java.lang.Object'
$wsClassFields.Range("D18").Value = 'java.lang.String'
$wsClassFields.Range("D19").Value = 'java.lang.String'
$wsClassFields.Range("D20").Value = 'java.lang.Object'
$wsClassFields.Range("D21").Value = 'java.lang.Long'
$wsClassFields.Range("D24").Value = 'java.lang.Integer'
$wsClassFields.Range("D27").Value = 'java.util.List'
$wsClassFields.Range("B37").Value = 'data'
$wsClassFields.Range("B38").Value = 'code'
$wsClassFields.Range("D38").Value = 'long'
$wsClassFields.Range("B39").Value = '$VALUES'
$wsClassFields.Range("D39").Value = 'com.macro.mall.common.api.ResultCode[]'
$wsClassFields.Range("B49").Value = 'contactName'
$wsClassFields.Range("B50").Value = 'enableSecurity'
$wsClassFields.Range("B51").Value = 'startTime'
$wsClassFields.Range("D51").Value = 'java.lang.String'
$wsClassFields.Range("D53").Value = 'boolean'
$wsClassFields.Range("B54").Value = 'version'
$wsClassFields.Range("B55").Value = 'apiBasePackage'
$wsClassFields.Range("B56").Value = 'contactUrl'
$wsClassFields.Range("B57").Value = 'enableSecurity'
$wsClassFields.Range("B58").Value = 'contactEmail'
$wsClassFields.Range("B59").Value = 'contactName'
$wsClassFields.Range("B60").Value = 'title'
$wsClassFields.Range("B61").Value = 'startTime'
$wsClassFields.Range("D63").Value = 'boolean'
$wsClassFields.Range("C64").Value = 'private'
$wsClassFields.Range("D64").Value = 'com.macro.mall.common.api.IErrorCode'
$wsClassFields.Range("C65").Value = ''
$wsClassFields.Range("D65").Value = 'long'
$wsClassFields.Range("B66").Value = 'errorCode'
$wsFieldClassRelations = $wb.Worksheets.Item("fieldClassRelations")
$wsFieldClassRelations.Range("D3").Value = 'FORBIDDEN'
$wsFieldClassRelations.Range("D4").Value = 'UNAUTHORIZED'
$wsFieldClassRelations.Range("D5").Value = 'VALIDATE_FAILED'
$wsFieldClassRelations.Range("D6").Value = 'FORBIDDEN'
$wsFieldClassRelations.Range("D7").Value = 'FAILED'
$wsFieldClassRelations.Range("D9").Value = 'UNAUTHORIZED'
